# Updated symbol list on Sat Dec 31 03:18:16 UTC 2022 with GitHub Actions
#
# Applies the refreshed coin price / volume snapshot to Sheet1. The "Price"
# column (D) holds numeric-looking values that must stay stored as TEXT
# (matching the original file's inlineStr cells), so each of those cells is
# forced to Text number-format before the value is written - otherwise Excel
# would silently reinterpret a string like "246.15" as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value2 = $value
}

function Set-PlainValue($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# --- Row 2 (BNB) ---
Set-TextValue 2 4 "246.15"

# --- Row 3 (OKB) ---
Set-TextValue 3 4 "25.34"

# --- Row 4 (HuobiToken) ---
Set-TextValue 4 4 "5.130"

# --- Row 5 (Cronos) ---
Set-TextValue 5 4 "0.05590"

# --- Row 6 (KuCoinToken) ---
Set-TextValue 6 4 "6.530"

# --- Row 7 (GateToken) ---
Set-TextValue 7 4 "3.017"

# --- Row 8 (MXToken) ---
Set-TextValue 8 4 "0.8174"

# --- Row 9 (FTXToken) ---
Set-TextValue 9 4 "0.8404"

# --- Row 10: was "One", now "WazirX" ---
Set-PlainValue 10 2 "WazirX"
Set-PlainValue 10 3 "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue  10 4 "0.1343"
Set-PlainValue 10 5 "9WazirXWRX"

# --- Row 11: was "WazirX", now "MandalaExchangeToken" ---
Set-PlainValue 11 2 "MandalaExchangeToken"
Set-PlainValue 11 3 "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue  11 4 "0.06946"
Set-PlainValue 11 5 "10MandalaExchangeTokenMDX"

# --- Row 12 (BitrueCoin) ---
Set-TextValue 12 4 "0.02846"

# --- Row 13 (BitMartToken) ---
Set-TextValue 13 4 "0.09388"

# --- Row 14 (BitForexToken) ---
Set-TextValue 14 4 "0.001529"

# --- Row 15: was "TigerCash", now "One" ---
Set-PlainValue 15 2 "One"
Set-PlainValue 15 3 "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue  15 4 "0.0005971"
Set-PlainValue 15 5 "14OneONEWorstin24h"

# --- Row 16: was "LEO", now "TigerCash" ---
Set-PlainValue 16 2 "TigerCash"
Set-PlainValue 16 3 "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue  16 4 "0.006148"
Set-PlainValue 16 5 "15TigerCashTCH"

# --- Row 17: was "BTSEToken", now "LEO" ---
Set-PlainValue 17 2 "LEO"
Set-PlainValue 17 3 "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue  17 4 "3.512"
Set-PlainValue 17 5 "16LEOLEO"

# --- Row 18: was "BitpandaEcosystemToken", now "BTSEToken" ---
Set-PlainValue 18 2 "BTSEToken"
Set-PlainValue 18 3 "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue  18 4 "2.082"
Set-PlainValue 18 5 "17BTSETokenBTSE"

# --- Row 19: was "MandalaExchangeToken", now "BitpandaEcosystemToken" ---
Set-PlainValue 19 2 "BitpandaEcosystemToken"
Set-PlainValue 19 3 "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue  19 4 "0.3179"
Set-PlainValue 19 5 "18BitpandaEcosystemTokenBEST"

# --- Row 20 (LiechtensteinCryptoassetsExchange) ---
Set-TextValue 20 4 "0.03167"

# --- Row 21 (ProBitToken) ---
Set-TextValue 21 4 "0.1329"

# --- Row 22 (MCDex) ---
Set-TextValue 22 4 "3.746"

# --- Row 23 (CoinExToken) ---
Set-TextValue 23 4 "0.04722"

# --- Row 25 (BitKan) ---
Set-TextValue 25 4 "0.001244"

# --- Row 26 (HotbitToken) ---
Set-TextValue 26 4 "0.004270"

# --- Row 27 (NitroEx) ---
Set-TextValue 27 4 "0.00009701"
Set-PlainValue 27 5 "26NitroExNTXBestin24h"

# --- Row 28 (UpBots) ---
Set-PlainValue 28 5 "27UpBotsUBXT"

# --- Row 40 (IDEX) ---
Set-TextValue 40 4 "0.03665"

# --- Row 41 (KickToken) ---
Set-TextValue 41 4 "0.006231"

# --- Row 42 (BKEXToken) ---
Set-TextValue 42 4 "0.1052"

# --- Row 43 (CEJI) ---
Set-TextValue 43 4 "0.002624"

# --- Row 44 (LocalTraders) ---
Set-TextValue 44 4 "0.008418"

# --- Row 45 (CoinLion) ---
Set-TextValue 45 4 "0.00005292"

# --- Row 48 (BOLO) ---
Set-TextValue 48 4 "0.002118"
